$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("2x Qty") -- shifts old D/E/F/G to E/F/G/H
$ws.Columns("D").Insert()

# Header for the new column
$ws.Range("D2").Value = "2x Qty"

# Fill the new "2x Qty" column with a doubling formula for every data row.
# Row 3 gets its own formula; rows 4-34 share the same formula pattern.
for ($r = 3; $r -le 34; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=C$r*2"
}

# --- Content corrections on the (now shifted) MFG# / DigiKey# columns (G/H) ---

# C10 C13 C14 C18 C19 C21 C22 C23 C4 C6 C7 C8 C28 C29 (100nF 0603) - row 6
$ws.Range("G6").Value = "C0603C104M5RACTU"
$ws.Range("H6").Value = "399-7845-1-ND"

# C20 C24 (22pF 0603) - row 7
$ws.Range("G7").Value = "C0603C220F1GACTU"
$ws.Range("H7").Value = "399-11145-1-ND"

# P1 (JTAGICE3 0.05") - row 14: MFG# becomes numeric, DigiKey# updated
$ws.Range("G14").Value = 62201021121
$ws.Range("H14").Value = "732-5374-ND"

# P3 (USB_OTG Custom) - row 16
$ws.Range("G16").Value = "10118192-0001LF"
$ws.Range("H16").Value = "609-4613-1-ND"

# --- View / window settings ---
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("B41").Select()
